# Quarterly financials update: insert a new "quarter ending 9/30/2018"
# column before column D, shifting the existing D:K data right to E:L,
# and populate the new column with the latest quarter figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VNRR")

# 1. Insert a new column at D; existing D:K shift to E:L automatically.
$ws.Columns.Item(4).Insert()

# 2. The freshly inserted column D has no number formatting - clone it
#    from column E (which now holds what used to be column D), so the
#    new cells pick up the same date / number styles as their row.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate column D with the new quarter's values.
function Set-D($row, $value) {
    $ws.Range("D$row").Value() = $value
}

# Income Statement (new quarter ending 2018-09-30)
Set-D 7  43373
Set-D 8  85500
Set-D 9  45200
Set-D 10 40300
Set-D 12 200
Set-D 13 0
Set-D 14 1000
Set-D 15 35600
Set-D 17 102200
Set-D 18 -16700
Set-D 20 700
Set-D 21 19500
Set-D 22 16100
Set-D 23 -32100
Set-D 24 0
Set-D 25 0
Set-D 26 -32100
Set-D 27 -32100
Set-D 28 0
Set-D 29 0
Set-D 30 0
Set-D 31 0
Set-D 32 -700
Set-D 33 -32100
Set-D 34 0
Set-D 35 -32100

# Balance Sheet
Set-D 38 43373
Set-D 41 4000
Set-D 42 0
Set-D 43 52100
Set-D 44 0
Set-D 45 40100
Set-D 46 96200
Set-D 47 0
Set-D 48 1399900
Set-D 49 "NA"
Set-D 50 0
Set-D 51 0
Set-D 52 11300
Set-D 53 0
Set-D 54 1507300
Set-D 57 44100
Set-D 58 "NA"
Set-D 59 144000
Set-D 60 188100
Set-D 61 863900
Set-D 62 181000
Set-D 63 0
Set-D 64 0
Set-D 65 0
Set-D 66 1233000
Set-D 68 0
Set-D 69 0
Set-D 70 0
Set-D 71 0
Set-D 72 -234000
Set-D 73 0
Set-D 74 0
Set-D 75 0
Set-D 76 274300
Set-D 77 0

# Cash Flow Statement
Set-D 80 43373
Set-D 81 -32100
Set-D 83 35600
Set-D 84 0
Set-D 85 0
Set-D 86 0
Set-D 87 0
Set-D 88 0
Set-D 89 11900
Set-D 91 -100
Set-D 92 0
Set-D 93 0
Set-D 94 12300
Set-D 96 0
Set-D 97 0
Set-D 98 0
Set-D 99 0
Set-D 100 -29500
Set-D 101 0
Set-D 102 -5300

Write-Output "done"
